$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2037.4
$ws.Range("I34").Value = 2171.75
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 2171.75
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -1968.75
$ws.Range("N34").Value = -1906
$ws.Range("H36").Value = 2037.4
$ws.Range("I36").Value = 2171.75
$ws.Range("J36").Value = 1500
$ws.Range("K36").Value = 2171.75
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = -1456.75
$ws.Range("N36").Value = -2930
$ws.Range("H88").Value = 4575.75
$ws.Range("I88").Value = 4399
$ws.Range("K88").Value = 4399
$ws.Range("M88").Value = -3993
$ws.Range("H91").Value = 4575.75
$ws.Range("I91").Value = 4399
$ws.Range("K91").Value = 4399
$ws.Range("M91").Value = -2995
$ws.Range("H111").Value = 3650
$ws.Range("I111").Value = 3737.5
$ws.Range("K111").Value = 11212.5
$ws.Range("M111").Value = -8145.5
$ws.Range("H132").Value = 3426.7368
$ws.Range("I132").Value = 1739.0769
$ws.Range("K132").Value = 5217.2307
$ws.Range("M132").Value = -2687.2307
$ws.Range("H138").Value = 8701.645500000001
$ws.Range("I138").Value = 7449
$ws.Range("J138").Value = 8990.718000000001
$ws.Range("K138").Value = 22347
$ws.Range("L138").Value = 26972.154
$ws.Range("M138").Value = -17207
$ws.Range("N138").Value = -37252.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 833.8333
$ws.Range("I2").Value = 833.8333
$ws.Range("K2").Value = 833.8333
$ws.Range("M2").Value = -720.8333
$ws.Range("H32").Value = 5327.3267
$ws.Range("I32").Value = 4245.3335
$ws.Range("K32").Value = 4245.3335
$ws.Range("M32").Value = -3958.3335
$ws.Range("H45").Value = 1822
$ws.Range("I45").Value = 1842.7142
$ws.Range("K45").Value = 1842.7142
$ws.Range("M45").Value = -1465.7142
$ws.Range("H74").Value = 2270.1428
$ws.Range("I74").Value = 2270.1428
$ws.Range("K74").Value = 2270.1428
$ws.Range("M74").Value = -1396.1428
$ws.Range("H77").Value = 2270.1428
$ws.Range("I77").Value = 2270.1428
$ws.Range("K77").Value = 11350.714
$ws.Range("M77").Value = -6982.714
$ws.Range("H116").Value = 833.8333
$ws.Range("I116").Value = 833.8333
$ws.Range("K116").Value = 833.8333
$ws.Range("M116").Value = 1460.1667
$ws.Range("H132").Value = 1706
$ws.Range("I132").Value = 1706
$ws.Range("K132").Value = 5118
$ws.Range("M132").Value = -2588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 833.8333
$ws.Range("I3").Value = 833.8333
$ws.Range("K3").Value = 833.8333
$ws.Range("M3").Value = -719.8333
$ws.Range("H80").Value = 790.0909
$ws.Range("I80").Value = 575.75
$ws.Range("J80").Value = 912.5714
$ws.Range("K80").Value = 575.75
$ws.Range("L80").Value = 912.5714
$ws.Range("M80").Value = 422.25
$ws.Range("N80").Value = -2908.5714
$ws.Range("H83").Value = 790.0909
$ws.Range("I83").Value = 575.75
$ws.Range("J83").Value = 912.5714
$ws.Range("K83").Value = 2878.75
$ws.Range("L83").Value = 4562.857
$ws.Range("M83").Value = 2113.25
$ws.Range("N83").Value = -14546.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 14522
$ws.Range("J2").Value = 16825.666
$ws.Range("L2").Value = 16825.666
$ws.Range("N2").Value = -17051.666
$ws.Range("H22").Value = 1632
$ws.Range("I22").Value = 1089.25
$ws.Range("K22").Value = 1089.25
$ws.Range("M22").Value = -739.25
$ws.Range("H31").Value = 2894.923
$ws.Range("I31").Value = 2785.0908
$ws.Range("J31").Value = 3499
$ws.Range("K31").Value = 2785.0908
$ws.Range("L31").Value = 3499
$ws.Range("M31").Value = -2490.0908
$ws.Range("N31").Value = -4089
$ws.Range("H34").Value = 2894.923
$ws.Range("I34").Value = 2785.0908
$ws.Range("J34").Value = 3499
$ws.Range("K34").Value = 2785.0908
$ws.Range("L34").Value = 3499
$ws.Range("M34").Value = -2583.0908
$ws.Range("N34").Value = -3903
$ws.Range("H94").Value = 3937.3
$ws.Range("I94").Value = 3851.2856
$ws.Range("J94").Value = 4138
$ws.Range("K94").Value = 3851.2856
$ws.Range("L94").Value = 4138
$ws.Range("M94").Value = -3400.2856
$ws.Range("N94").Value = -5040
$ws.Range("H122").Value = 938.1818
$ws.Range("I122").Value = 1035.5555
$ws.Range("K122").Value = 3106.6665
$ws.Range("M122").Value = -656.6664999999998
$ws.Range("H132").Value = 10545.588
$ws.Range("I132").Value = 4182.2856
$ws.Range("K132").Value = 12546.8568
$ws.Range("M132").Value = -10016.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1876.4286
$ws.Range("I113").Value = 1439.1666
$ws.Range("K113").Value = 1439.1666
$ws.Range("M113").Value = 730.8334
$ws.Range("H122").Value = 3502.2727
$ws.Range("I122").Value = 3369.8235
$ws.Range("J122").Value = 3952.6
$ws.Range("K122").Value = 10109.4705
$ws.Range("L122").Value = 11857.8
$ws.Range("M122").Value = -7659.470499999999
$ws.Range("N122").Value = -16757.8
$ws.Range("H126").Value = 1000000000
$ws.Range("I126").Value = 1000000000
$ws.Range("K126").Value = 3000000000
$ws.Range("M126").Value = -2999997530
$ws.Range("H132").Value = 3968.4
$ws.Range("I132").Value = 3968.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11905.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9375.200000000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 14999.5
$ws.Range("J3").Value = 14999.5
$ws.Range("L3").Value = 14999.5
$ws.Range("N3").Value = -15223.5
$ws.Range("H15").Value = 14999.5
$ws.Range("J15").Value = 14999.5
$ws.Range("L15").Value = 14999.5
$ws.Range("N15").Value = -15339.5
$ws.Range("H16").Value = 1036.7646
$ws.Range("J16").Value = 833
$ws.Range("L16").Value = 833
$ws.Range("N16").Value = -1173
$ws.Range("H82").Value = 3081.2083
$ws.Range("J82").Value = 3669.2307
$ws.Range("L82").Value = 3669.2307
$ws.Range("N82").Value = -4391.2307
$ws.Range("H85").Value = 3081.2083
$ws.Range("J85").Value = 3669.2307
$ws.Range("L85").Value = 3669.2307
$ws.Range("N85").Value = -6165.2307
$ws.Range("H132").Value = 4854
$ws.Range("I132").Value = 4508.5
$ws.Range("K132").Value = 13525.5
$ws.Range("M132").Value = -10995.5
$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -55059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1890.0834
$ws.Range("I126").Value = 1497.4286
$ws.Range("K126").Value = 4492.2858
$ws.Range("M126").Value = -2022.2858
$ws.Range("H136").Value = 2239.2917
$ws.Range("I136").Value = 1684.4783
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 5053.4349
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -2503.4349
$ws.Range("N136").Value = -50100
